# Fix the "Beneign" -> "Benign" classification typo in the
# publication_parameters sheet (column G holds the classification label).
# Every cell that used the misspelled "Beneign" string gets corrected to
# "Benign"; Excel will add the new, correctly spelled value to the shared
# strings table automatically the first time it's written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("publication_parameters")

$fixedCells = @("G20","G21","G22","G23","G24","G25","G26","G27","G28","G30","G47")
foreach ($addr in $fixedCells) {
    $ws.Range($addr).Value = "Benign"
}

# Leave the sheet on the cell the author ended up selecting after editing.
$ws.Activate()
$ws.Range("E67").Select()
